{"js": "// Commit message: \"Work on chapter 2\"\n//\n// The underlying change is Word's automatic relocation of the `_GoBack`\n// bookmark (the \"last edit position\" marker) from the very start of the\n// document (next to the \"Dissertation Title\" heading) to the end of the\n// paragraph that was most recently edited \u2014 the paragraph ending in\n// \"...upstream of their target gene\" in the \"Chapter X: Groucho activity\n// in the developing embryo\" > Results section. All other bookmark ids\n// shift accordingly (Word renumbers bookmark ids sequentially by their\n// position in the document whenever the set changes).\n\n// 1) Remove the old `_GoBack` bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Find the paragraph that ends with the \"...upstream of their target\n//    gene\" sentence (the last-edited location) so we can drop the\n//    `_GoBack` bookmark back in at its end.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text && paragraph.text.indexOf(\"upstream of their target gene\") !== -1) {\n    targetParagraph = paragraph;\n    break;\n  }\n}\n\n// 3) Re-insert `_GoBack` collapsed at the end of that paragraph's content.\nif (targetParagraph) {\n  const endRange = targetParagraph.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Commit message: \"Work on chapter 2\"\n#\n# Word stores the position of the user's last edit in a hidden bookmark\n# named \"_GoBack\". This change simply reflects that the last edit moved\n# from the very top of the document (next to the \"Dissertation Title\"\n# heading) down into Chapter 2 / Results, specifically to the end of the\n# paragraph ending \"...upstream of their target gene\". Re-saving the\n# document relocates the \"_GoBack\" bookmark accordingly (and, because\n# bookmark ids are assigned sequentially by document order, every other\n# bookmark whose id came after the old \"_GoBack\" position shifts down by\n# one to fill the gap).\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark, wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the paragraph that was last edited - the one ending in\n#    \"...upstream of their target gene\".\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"upstream of their target gene\")\nif ($found) {\n    $targetParagraph = $searchRange.Paragraphs(1)\n    $paragraphEnd = $targetParagraph.Range.End\n\n    # 3) Insert a short, unique marker right at the end of the paragraph's\n    #    text (i.e. immediately before its paragraph mark). Placing a\n    #    brand-new, *collapsed* bookmark exactly at \"paragraph end minus\n    #    one\" is unreliable, so we temporarily insert text there, bookmark\n    #    in front of that text (a non-boundary position), then delete the\n    #    marker text again - the bookmark itself stays collapsed in place.\n    $marker = \"@@GOBACK_MARKER@@\"\n    $insertionPoint = $d.Range($paragraphEnd - 1, $paragraphEnd - 1)\n    $insertionPoint.InsertAfter($marker)\n\n    # 4) Find the marker and add the \"_GoBack\" bookmark immediately before it.\n    $markerRange = $d.Content\n    $markerRange.Find.Execute($marker) | Out-Null\n    $bookmarkPosition = $d.Range($markerRange.Start, $markerRange.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkPosition)\n\n    # 5) Remove the temporary marker text; the bookmark remains collapsed\n    #    at the end of the paragraph's real content.\n    $markerRange2 = $d.Content\n    $markerRange2.Find.Execute($marker) | Out-Null\n    $markerRange2.Delete()\n}\n"}
